$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.150782
$ws.Range("H2").Value = 3.452345999999999
$ws.Range("I2").Value = 0.03823856951930295
$ws.Range("J2").Value = 0.03823856951930295
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 193.4577633183079
$ws.Range("R2").Value = 1741.119869864772
$ws.Range("S2").Value = 0.0114110860732047
$ws.Range("T2").Value = 0.0114110860732047

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.150782
$ws.Range("H3").Value = 3.452345999999999
$ws.Range("I3").Value = 0.03823856951930295
$ws.Range("J3").Value = 0.03823856951930295
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 187.584641892958
$ws.Range("R3").Value = 1688.261777036622
$ws.Range("S3").Value = 0.01106466061602218
$ws.Range("T3").Value = 0.01106466061602218

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.150782
$ws.Range("H4").Value = 3.452345999999999
$ws.Range("I4").Value = 0.03823856951930295
$ws.Range("J4").Value = 0.03823856951930295
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 191.02236644046
$ws.Range("R4").Value = 1719.20129796414
$ws.Range("S4").Value = 0.0112674344413505
$ws.Range("T4").Value = 0.0112674344413505

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.150782
$ws.Range("H5").Value = 3.452345999999999
$ws.Range("I5").Value = 0.03823856951930295
$ws.Range("J5").Value = 0.03823856951930295
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 76.212533789582
$ws.Range("R5").Value = 685.9128041062379
$ws.Range("S5").Value = 0.004495388388725572
$ws.Range("T5").Value = 0.004495388388725571

$ws.Range("I6").Value = 0.9169230158851821
$ws.Range("J6").Value = 0.916923015885182
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 4638.925514686937
$ws.Range("R6").Value = 41750.32963218244
$ws.Range("S6").Value = 0.273626539598623
$ws.Range("T6").Value = 0.273626539598623

$ws.Range("I7").Value = 0.9169230158851821
$ws.Range("J7").Value = 0.916923015885182
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.2653196003231137
$ws.Range("T7").Value = 0.2653196003231137

$ws.Range("I8").Value = 0.9169230158851821
$ws.Range("J8").Value = 0.916923015885182
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 4580.527110191537
$ws.Range("R8").Value = 41224.74399172383
$ws.Range("S8").Value = 0.2701819157758076
$ws.Range("T8").Value = 0.2701819157758076

$ws.Range("I9").Value = 0.9169230158851821
$ws.Range("J9").Value = 0.916923015885182
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 1827.50106003099
$ws.Range("R9").Value = 16447.50954027891
$ws.Range("S9").Value = 0.1077949601876378
$ws.Range("T9").Value = 0.1077949601876377

$ws.Range("G10").Value = 1.290098666666667
$ws.Range("H10").Value = 3.870296
$ws.Range("I10").Value = 0.04286783035543951
$ws.Range("J10").Value = 0.0428678303554395
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 216.8782640962968
$ws.Range("R10").Value = 1951.904376866672
$ws.Range("S10").Value = 0.01279254187870505
$ws.Range("T10").Value = 0.01279254187870505

$ws.Range("G11").Value = 1.290098666666667
$ws.Range("H11").Value = 3.870296
$ws.Range("I11").Value = 0.04286783035543951
$ws.Range("J11").Value = 0.0428678303554395
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 210.2941272919191
$ws.Range("R11").Value = 1892.647145627272
$ws.Range("S11").Value = 0.01240417725324988
$ws.Range("T11").Value = 0.01240417725324987

$ws.Range("G12").Value = 1.290098666666667
$ws.Range("H12").Value = 3.870296
$ws.Range("I12").Value = 0.04286783035543951
$ws.Range("J12").Value = 0.0428678303554395
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 214.1480317282933
$ws.Range("R12").Value = 1927.33228555464
$ws.Range("S12").Value = 0.01263149940609113
$ws.Range("T12").Value = 0.01263149940609113

$ws.Range("G13").Value = 1.290098666666667
$ws.Range("H13").Value = 3.870296
$ws.Range("I13").Value = 0.04286783035543951
$ws.Range("J13").Value = 0.0428678303554395
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 85.43902166112089
$ws.Range("R13").Value = 768.9511949500879
$ws.Range("S13").Value = 0.005039611817393456
$ws.Range("T13").Value = 0.005039611817393455

$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.05930433333333333
$ws.Range("H14").Value = 0.177913
$ws.Range("I14").Value = 0.001970584240075516
$ws.Range("J14").Value = 0.001970584240075516
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 9.96964123678511
$ws.Range("R14").Value = 89.72677113106599
$ws.Range("S14").Value = 0.0005880582527191849
$ws.Range("T14").Value = 0.0005880582527191849

$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.05930433333333333
$ws.Range("H15").Value = 0.177913
$ws.Range("I15").Value = 0.001970584240075516
$ws.Range("J15").Value = 0.001970584240075516
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 9.666976135387888
$ws.Range("R15").Value = 87.002785218491
$ws.Range("S15").Value = 0.0005702055831537033
$ws.Range("T15").Value = 0.0005702055831537032

$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.05930433333333333
$ws.Range("H16").Value = 0.177913
$ws.Range("I16").Value = 0.001970584240075516
$ws.Range("J16").Value = 0.001970584240075516
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 9.844135634296665
$ws.Range("R16").Value = 88.59722070866999
$ws.Range("S16").Value = 0.000580655317793753
$ws.Range("T16").Value = 0.0005806553177937529

$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.05930433333333333
$ws.Range("H17").Value = 0.177913
$ws.Range("I17").Value = 0.001970584240075516
$ws.Range("J17").Value = 0.001970584240075516
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 3.927532328482111
$ws.Range("R17").Value = 35.347790956339
$ws.Range("S17").Value = 0.0002316650864088747
$ws.Range("T17").Value = 0.0002316650864088746
